$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 07:07:00"
$wsZhCn.Range("H2").Value = "2016-03-13 07:07:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 07:07:04"
$wsDeDe.Range("H2").Value = "2016-03-13 07:07:24"
